# "Fruta / hortaliza, semanal" weekly data refresh.
# The weekly pull updates the Fecha (D), Volumen (J), Precio mínimo (K),
# Precio máximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# columns for each data row (rows 2-9). Other columns (market, region,
# category, quality, unit, origin, classification, etc.) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$rows = @(
    @{ Row = 2; D = 45092; J = 210; K = 10000; L = 11000; M = 10714; P = 595 },
    @{ Row = 3; D = 45245; J = 100; K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 4; D = 44792; J = 160; K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 5; D = 45205; J = 200; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 6; D = 44804; J = 50;  K = 9500;  L = 10000; M = 9750;  P = 542 },
    @{ Row = 7; D = 45215; J = 200; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 8; D = 45175; J = 250; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 9; D = 44714; J = 80;  K = 9000;  L = 10000; M = 9500;  P = 528 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("P$n").Value = $r.P
}
